$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add New York state hospitalization data for 22 April 2020 (row 39).
# Carry down the formatting of the row above (row 38) cell-by-cell so the
# new row matches the existing table style: column A keeps the date
# number format, the other populated columns keep the default (General)
# numeric style used throughout the sheet.
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("B38").Copy()
$ws.Range("B39").PasteSpecial(-4122)
$ws.Range("D38").Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("F38").Copy()
$ws.Range("F39").PasteSpecial(-4122)
$ws.Range("G38").Copy()
$ws.Range("G39").PasteSpecial(-4122)

$ws.Range("A39").Value = 43943
$ws.Range("B39").Value = -578
$ws.Range("D39").Value = -16
$ws.Range("F39").Value = 438
$ws.Range("G39").Value = 1359

# Move the active selection down one row, matching the sheet's prior
# "next empty row" selection convention (F39 -> F40)
$ws.Range("F40").Select()
